$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.244.50"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.789.72"

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'225.86"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D5").ClearFormats()

$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'32.28"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D8").ClearFormats()

$ws.Range("D9").Value = "'0.295"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D9").ClearFormats()

$ws.Range("D10").Value = "'0.0689"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D10").ClearFormats()

$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D11").ClearFormats()

$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.00"
$ws.Range("E13").Value = "  -4.65%  "
$ws.Range("D13").ClearFormats()

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.778.70"
$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").Value = "'0.627"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = "34.206.68"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "'67.94"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  +2.36%  "

$ws.Range("D20").Value = "'246.56"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D20").ClearFormats()

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "'4.18"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D23").ClearFormats()

$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D24").ClearFormats()

$ws.Range("D25").Value = "'162.28"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D25").ClearFormats()

$ws.Range("D26").Value = "'7.18"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D26").ClearFormats()

$ws.Range("D27").Value = "'16.36"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D27").ClearFormats()

$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  -0.62%  "

$ws.Range("D31").Value = "'0.0521"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D31").ClearFormats()

$ws.Range("D32").Value = "'3.76"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D32").ClearFormats()

$ws.Range("D33").Value = "'3.82"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D33").ClearFormats()

$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("D35").Value = "1.442.00"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("E36").Value = "  +9.80%  "

$ws.Range("E37").Value = "  +1.73%  "

$ws.Range("E38").Value = "  +1.60%  "

$ws.Range("D40").Value = "'82.27"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D40").ClearFormats()

$ws.Range("E41").Value = "  +2.00%  "

$ws.Range("D42").Value = "'14.11"
$ws.Range("E42").Value = "  +5.75%  "
$ws.Range("D42").ClearFormats()

$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D43").ClearFormats()

$ws.Range("D44").Value = "'0.923"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D44").ClearFormats()

$ws.Range("D45").Value = "'0.0519"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D45").ClearFormats()

$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").Value = "1.940.67"
$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("D49").Value = "'105.62"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D49").ClearFormats()

$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("E51").Value = "  -6.47%  "
